$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Rol" column header in W1, styled like the existing pink P/H1/H2-3
# headers (T1) but with only left+right borders (no top/bottom).
$src = $ws.Range("T1")
$dst = $ws.Range("W1")
$src.Copy()
$dst.PasteSpecial(-4122)
$dst.Value = "Rol"
$dst.Borders.Item(8).LineStyle = -4142
$dst.Borders.Item(9).LineStyle = -4142

# Fill the new column with "Studenten" for every data row.
$ws.Range("W2").Value = "Studenten"
$ws.Range("W3").Value = "Studenten"
$ws.Range("W4").Value = "Studenten"
$ws.Range("W5").Value = "Studenten"

# Reflect the user's final selection/view after adding the column.
$ws.Range("W2:W5").Select()

Write-Host "done"
